$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on changed cells so numeric-looking strings
# (e.g. "0.0868", "63.46") are preserved exactly as text, matching
# the inlineStr cell type used throughout column D and E.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.702.23'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.589.34'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.40%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.60%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.31'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.08%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.25'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -4.05%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.78%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0868'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.815.40'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.620.51'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.83%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.692.80'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.46'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '219.88'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0695'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.78%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.55%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.27%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.57'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.75%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.15'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.33%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0468'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.23'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.370.52'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.86%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.51%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.978'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.94%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0167'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.54%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.71%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.974'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '64.27'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.54%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.29%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.726.18'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.80'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0100'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +11.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0969'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0496'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.37%  '
